# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "last updated" timestamp string in A1 (23:22 -> 23:52)
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 23:52"

# 2. Insert "Gabon" as a new row right after "Martinica" (row 130) and before
#    "Ruanda" (previously row 131). This shifts Ruanda..Madagascar down by one
#    row (rows 131-138), and the old "Gabon" row (previously row 138, just
#    before "Trinidad yTobago") disappears since Gabon now lives at row 131.
#    Row 139 ("Trinidad yTobago") and beyond are untouched.

# Row 131: Gabon (new data)
$ws.Range("A131").Value = "Gabon"
$ws.Range("B131").Value = 156
$ws.Range("C131").Value = 36
$ws.Range("D131").Value = 16
$ws.Range("E131").Value = 139
$ws.Range("F131").Value = 2
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 1

# Row 132: Ruanda (was row 131's data)
$ws.Range("A132").Value = "Ruanda"
$ws.Range("B132").Value = 150
$ws.Range("C132").Value = 3
$ws.Range("D132").Value = 84
$ws.Range("E132").Value = 66
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 0

# Row 133: Guadalupe (was row 132's data)
$ws.Range("A133").Value = "Guadalupe"
$ws.Range("B133").Value = 148
$ws.Range("C133").Value = 0
$ws.Range("D133").Value = 73
$ws.Range("E133").Value = 63
$ws.Range("F133").Value = 13
$ws.Range("G133").Value = 4
$ws.Range("H133").Value = 12

# Row 134: Brunei (was row 133's data)
$ws.Range("A134").Value = "Brunei"
$ws.Range("B134").Value = 138
$ws.Range("C134").Value = 0
$ws.Range("D134").Value = 116
$ws.Range("E134").Value = 21
$ws.Range("F134").Value = 2
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 1

# Row 135: Gibraltar (was row 134's data)
$ws.Range("A135").Value = "Gibraltar"
$ws.Range("B135").Value = 132
$ws.Range("C135").Value = 0
$ws.Range("D135").Value = 120
$ws.Range("E135").Value = 12
$ws.Range("F135").Value = 1
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 0

# Row 136: Camboya (was row 135's data)
$ws.Range("A136").Value = "Camboya"
$ws.Range("B136").Value = 122
$ws.Range("C136").Value = 0
$ws.Range("D136").Value = 110
$ws.Range("E136").Value = 12
$ws.Range("F136").Value = 1
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 0

# Row 137: Birmania (was row 136's data)
$ws.Range("A137").Value = "Birmania"
$ws.Range("B137").Value = 121
$ws.Range("C137").Value = 2
$ws.Range("D137").Value = 7
$ws.Range("E137").Value = 109
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 5

# Row 138: Madagascar (was row 137's data)
$ws.Range("A138").Value = "Madagascar"
$ws.Range("B138").Value = 121
$ws.Range("C138").Value = 0
$ws.Range("D138").Value = 44
$ws.Range("E138").Value = 77
$ws.Range("F138").Value = 1
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 0

# 3. Update United States row (row 4) stats
$ws.Range("B4").Value = 815892
$ws.Range("C4").Value = 23133
$ws.Range("E4").Value = 688146
$ws.Range("G4").Value = 2612
$ws.Range("H4").Value = 45126

# 4. Update row 20 stats
$ws.Range("B20").Value = 20080
$ws.Range("C20").Value = 1541
$ws.Range("D20").Value = 3975
$ws.Range("E20").Value = 15460
$ws.Range("G20").Value = 53
$ws.Range("H20").Value = 645
